# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off (but the handback is stale), refreshing the Overview sheet
# and the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77fa9888f2ef475e4b96c658fee4001221fb6ac3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5900506a846aa3730d4aa9cab16fcc301d8a35e/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 12:37:33"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-19 12:37:29"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P:P").ColumnWidth = 39.1

# ---------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-19 12:37:33"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P:P").ColumnWidth = 39.1
